$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows whose match data was reshuffled (F:V columns; A and E stay as-is) ---
# Row 13
$ws.Range("F13").Value = "Horsens"
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = "Koge"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1.72
$ws.Range("K13").Value = "01/08/2023 05:42"
$ws.Range("L13").Value = 1.81
$ws.Range("M13").Value = "04/08/2023 18:52"
$ws.Range("N13").Value = 3.95
$ws.Range("O13").Value = "01/08/2023 05:42"
$ws.Range("P13").Value = 4.16
$ws.Range("Q13").Value = "04/08/2023 18:56"
$ws.Range("R13").Value = 4.12
$ws.Range("S13").Value = "01/08/2023 05:42"
$ws.Range("T13").Value = 3.96
$ws.Range("U13").Value = "04/08/2023 18:56"
$ws.Range("V13").Value = "https://www.betexplorer.com/football/denmark/1st-division/horsens-koge/zDaoiXfI/"

# Row 14
$ws.Range("F14").Value = "Hillerod"
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = "Hobro"
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 2.23
$ws.Range("K14").Value = "29/07/2023 14:12"
$ws.Range("L14").Value = 3.15
$ws.Range("M14").Value = "04/08/2023 18:23"
$ws.Range("N14").Value = 3.38
$ws.Range("O14").Value = "29/07/2023 14:12"
$ws.Range("P14").Value = 3.49
$ws.Range("Q14").Value = "04/08/2023 18:23"
$ws.Range("R14").Value = 3.23
$ws.Range("S14").Value = "29/07/2023 14:12"
$ws.Range("T14").Value = 2.27
$ws.Range("U14").Value = "04/08/2023 18:23"
$ws.Range("V14").Value = "https://www.betexplorer.com/football/denmark/1st-division/hillerod-hobro/E3bkjiAO/"

# Row 26
$ws.Range("F26").Value = "Kolding IF"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = "Sonderjyske"
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 2.08
$ws.Range("K26").Value = "14/08/2023 04:12"
$ws.Range("L26").Value = 2.99
$ws.Range("M26").Value = "18/08/2023 18:30"
$ws.Range("N26").Value = 3.6
$ws.Range("O26").Value = "14/08/2023 04:12"
$ws.Range("P26").Value = 3.45
$ws.Range("Q26").Value = "18/08/2023 18:31"
$ws.Range("R26").Value = 3.38
$ws.Range("S26").Value = "14/08/2023 04:12"
$ws.Range("T26").Value = 2.38
$ws.Range("U26").Value = "18/08/2023 18:30"
$ws.Range("V26").Value = "https://www.betexplorer.com/football/denmark/1st-division/kolding-if-sonderjyske/UiWf7qoG/"

# Row 27
$ws.Range("F27").Value = "B.93"
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = "Fredericia"
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 2.04
$ws.Range("K27").Value = "14/08/2023 04:12"
$ws.Range("L27").Value = 2.44
$ws.Range("M27").Value = "18/08/2023 18:47"
$ws.Range("N27").Value = 3.66
$ws.Range("O27").Value = "14/08/2023 04:12"
$ws.Range("P27").Value = 3.73
$ws.Range("Q27").Value = "18/08/2023 18:47"
$ws.Range("R27").Value = 3.43
$ws.Range("S27").Value = "14/08/2023 04:12"
$ws.Range("T27").Value = 2.73
$ws.Range("U27").Value = "18/08/2023 18:37"
$ws.Range("V27").Value = "https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-fredericia/8xrn95G3/"

# Row 28
$ws.Range("F28").Value = "Horsens"
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = "Helsingor"
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = 2.01
$ws.Range("K28").Value = "13/08/2023 22:12"
$ws.Range("L28").Value = 2.31
$ws.Range("M28").Value = "18/08/2023 18:51"
$ws.Range("N28").Value = 3.79
$ws.Range("O28").Value = "13/08/2023 22:12"
$ws.Range("P28").Value = 3.76
$ws.Range("Q28").Value = "18/08/2023 18:51"
$ws.Range("R28").Value = 3.2
$ws.Range("S28").Value = "13/08/2023 22:12"
$ws.Range("T28").Value = 2.89
$ws.Range("U28").Value = "18/08/2023 18:51"
$ws.Range("V28").Value = "https://www.betexplorer.com/football/denmark/1st-division/horsens-helsingor/nTtj8PV9/"

# Row 31
$ws.Range("F31").Value = "Naestved"
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = "B.93"
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = 2.18
$ws.Range("K31").Value = "20/08/2023 15:13"
$ws.Range("L31").Value = 2.08
$ws.Range("M31").Value = "22/08/2023 18:56"
$ws.Range("N31").Value = 3.57
$ws.Range("O31").Value = "20/08/2023 15:13"
$ws.Range("P31").Value = 3.97
$ws.Range("Q31").Value = "22/08/2023 18:56"
$ws.Range("R31").Value = 3.16
$ws.Range("S31").Value = "20/08/2023 15:13"
$ws.Range("T31").Value = 3.2
$ws.Range("U31").Value = "22/08/2023 18:56"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/denmark/1st-division/naestved-if-boldklubben-1893/KxwSa1wj/"

# Row 32
$ws.Range("F32").Value = "Fredericia"
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = "Horsens"
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2.6
$ws.Range("K32").Value = "20/08/2023 15:13"
$ws.Range("L32").Value = 2.21
$ws.Range("M32").Value = "22/08/2023 18:57"
$ws.Range("N32").Value = 3.43
$ws.Range("O32").Value = "20/08/2023 15:13"
$ws.Range("P32").Value = 3.87
$ws.Range("Q32").Value = "22/08/2023 18:51"
$ws.Range("R32").Value = 2.65
$ws.Range("S32").Value = "20/08/2023 15:13"
$ws.Range("T32").Value = 3
$ws.Range("U32").Value = "22/08/2023 18:57"
$ws.Range("V32").Value = "https://www.betexplorer.com/football/denmark/1st-division/fredericia-horsens/vRLU0sOq/"

# Row 33
$ws.Range("F33").Value = "Sonderjyske"
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = "Koge"
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 1.4
$ws.Range("K33").Value = "20/08/2023 13:12"
$ws.Range("L33").Value = 1.36
$ws.Range("M33").Value = "23/08/2023 18:52"
$ws.Range("N33").Value = 4.75
$ws.Range("O33").Value = "20/08/2023 13:12"
$ws.Range("P33").Value = 5.5
$ws.Range("Q33").Value = "23/08/2023 18:55"
$ws.Range("R33").Value = 6.3
$ws.Range("S33").Value = "20/08/2023 13:12"
$ws.Range("T33").Value = 7.36
$ws.Range("U33").Value = "23/08/2023 18:55"
$ws.Range("V33").Value = "https://www.betexplorer.com/football/denmark/1st-division/sonderjyske-koge/QZTHMOp3/"

# Row 34
$ws.Range("F34").Value = "Helsingor"
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = "Hillerod"
$ws.Range("I34").Value = 6
$ws.Range("J34").Value = 2.12
$ws.Range("K34").Value = "19/08/2023 13:12"
$ws.Range("L34").Value = 1.82
$ws.Range("M34").Value = "23/08/2023 18:52"
$ws.Range("N34").Value = 3.52
$ws.Range("O34").Value = "19/08/2023 13:12"
$ws.Range("P34").Value = 3.97
$ws.Range("Q34").Value = "23/08/2023 18:52"
$ws.Range("R34").Value = 3.35
$ws.Range("S34").Value = "19/08/2023 13:12"
$ws.Range("T34").Value = 4.08
$ws.Range("U34").Value = "23/08/2023 18:52"
$ws.Range("V34").Value = "https://www.betexplorer.com/football/denmark/1st-division/helsingor-hillerod/bTyWbLhd/"

# Row 35
$ws.Range("F35").Value = "Vendsyssel"
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = "Kolding IF"
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = 2.19
$ws.Range("K35").Value = "20/08/2023 15:13"
$ws.Range("L35").Value = 2.73
$ws.Range("M35").Value = "23/08/2023 18:47"
$ws.Range("N35").Value = 3.53
$ws.Range("O35").Value = "20/08/2023 15:13"
$ws.Range("P35").Value = 3.4
$ws.Range("Q35").Value = "23/08/2023 18:47"
$ws.Range("R35").Value = 3.18
$ws.Range("S35").Value = "20/08/2023 15:13"
$ws.Range("T35").Value = 2.61
$ws.Range("U35").Value = "23/08/2023 18:08"
$ws.Range("V35").Value = "https://www.betexplorer.com/football/denmark/1st-division/vendsyssel-ff-kolding-if/CWPDN4Vd/"

# Row 43
$ws.Range("F43").Value = "Hillerod"
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = "Sonderjyske"
$ws.Range("I43").Value = 2
$ws.Range("J43").Value = 3.9
$ws.Range("K43").Value = "28/08/2023 18:42"
$ws.Range("L43").Value = 4.16
$ws.Range("M43").Value = "01/09/2023 18:58"
$ws.Range("N43").Value = 3.8
$ws.Range("O43").Value = "28/08/2023 18:42"
$ws.Range("P43").Value = 3.85
$ws.Range("Q43").Value = "01/09/2023 18:58"
$ws.Range("R43").Value = 1.79
$ws.Range("S43").Value = "28/08/2023 18:42"
$ws.Range("T43").Value = 1.83
$ws.Range("U43").Value = "01/09/2023 18:58"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/denmark/1st-division/hillerod-sonderjyske/EwHtnuEE/"

# Row 44
$ws.Range("F44").Value = "Vendsyssel"
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = "Fredericia"
$ws.Range("I44").Value = 2
$ws.Range("J44").Value = 2.02
$ws.Range("K44").Value = "27/08/2023 13:12"
$ws.Range("L44").Value = 2.46
$ws.Range("M44").Value = "01/09/2023 18:55"
$ws.Range("N44").Value = 3.79
$ws.Range("O44").Value = "27/08/2023 13:12"
$ws.Range("P44").Value = 3.9
$ws.Range("Q44").Value = "01/09/2023 18:55"
$ws.Range("R44").Value = 3.19
$ws.Range("S44").Value = "27/08/2023 13:12"
$ws.Range("T44").Value = 2.62
$ws.Range("U44").Value = "01/09/2023 18:55"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/denmark/1st-division/vendsyssel-ff-fredericia/zmIxmLb8/"

# Row 45
$ws.Range("F45").Value = "Horsens"
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = "B.93"
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1.79
$ws.Range("K45").Value = "28/08/2023 01:42"
$ws.Range("L45").Value = 1.59
$ws.Range("M45").Value = "01/09/2023 18:57"
$ws.Range("N45").Value = 4.02
$ws.Range("O45").Value = "28/08/2023 01:42"
$ws.Range("P45").Value = 4.54
$ws.Range("Q45").Value = "01/09/2023 18:57"
$ws.Range("R45").Value = 3.72
$ws.Range("S45").Value = "28/08/2023 01:42"
$ws.Range("T45").Value = 4.97
$ws.Range("U45").Value = "01/09/2023 18:57"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/denmark/1st-division/horsens-boldklubben-1893/MqUfOyM7/"

# --- Append new rows 55-58 (copy formatting from row 54, then set values) ---
# Row 55
$ws.Range("A54:V54").Copy()
$ws.Range("A55:V55").PasteSpecial(-4122)
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "denmark"
$ws.Range("C55").Value = "1st-division"
$ws.Range("D55").Value = "2023-2024"
$ws.Range("E55").Value = 45191.79166666666
$ws.Range("F55").Value = "Aalborg"
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = "Koge"
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1.18
$ws.Range("K55").Value = "20/09/2023 00:42"
$ws.Range("L55").Value = 1.21
$ws.Range("M55").Value = "22/09/2023 18:51"
$ws.Range("N55").Value = 6.97
$ws.Range("O55").Value = "20/09/2023 00:42"
$ws.Range("P55").Value = 7.43
$ws.Range("Q55").Value = "22/09/2023 18:51"
$ws.Range("R55").Value = 13.84
$ws.Range("S55").Value = "20/09/2023 00:42"
$ws.Range("T55").Value = 10.88
$ws.Range("U55").Value = "22/09/2023 18:51"
$ws.Range("V55").Value = "https://www.betexplorer.com/football/denmark/1st-division/aalborg-koge/preqskkM/"

# Row 56
$ws.Range("A54:V54").Copy()
$ws.Range("A56:V56").PasteSpecial(-4122)
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = "denmark"
$ws.Range("C56").Value = "1st-division"
$ws.Range("D56").Value = "2023-2024"
$ws.Range("E56").Value = 45191.79166666666
$ws.Range("F56").Value = "B.93"
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = "Hobro"
$ws.Range("I56").Value = 2
$ws.Range("J56").Value = 2.74
$ws.Range("K56").Value = "19/09/2023 23:42"
$ws.Range("L56").Value = 2.92
$ws.Range("M56").Value = "22/09/2023 18:54"
$ws.Range("N56").Value = 3.45
$ws.Range("O56").Value = "19/09/2023 23:42"
$ws.Range("P56").Value = 3.46
$ws.Range("Q56").Value = "22/09/2023 18:54"
$ws.Range("R56").Value = 2.41
$ws.Range("S56").Value = "19/09/2023 23:42"
$ws.Range("T56").Value = 2.42
$ws.Range("U56").Value = "22/09/2023 18:54"
$ws.Range("V56").Value = "https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-hobro/boamt94S/"

# Row 57
$ws.Range("A54:V54").Copy()
$ws.Range("A57:V57").PasteSpecial(-4122)
$ws.Range("A57").Value = 56
$ws.Range("B57").Value = "denmark"
$ws.Range("C57").Value = "1st-division"
$ws.Range("D57").Value = "2023-2024"
$ws.Range("E57").Value = 45191.79166666666
$ws.Range("F57").Value = "Kolding IF"
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = "Fredericia"
$ws.Range("I57").Value = 2
$ws.Range("J57").Value = 2.03
$ws.Range("K57").Value = "20/09/2023 00:42"
$ws.Range("L57").Value = 2.44
$ws.Range("M57").Value = "22/09/2023 18:37"
$ws.Range("N57").Value = 3.62
$ws.Range("O57").Value = "20/09/2023 00:42"
$ws.Range("P57").Value = 3.7
$ws.Range("Q57").Value = "22/09/2023 18:37"
$ws.Range("R57").Value = 3.5
$ws.Range("S57").Value = "20/09/2023 00:42"
$ws.Range("T57").Value = 2.74
$ws.Range("U57").Value = "22/09/2023 18:37"
$ws.Range("V57").Value = "https://www.betexplorer.com/football/denmark/1st-division/kolding-if-fredericia/4SyGyTCq/"

# Row 58
$ws.Range("A54:V54").Copy()
$ws.Range("A58:V58").PasteSpecial(-4122)
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "denmark"
$ws.Range("C58").Value = "1st-division"
$ws.Range("D58").Value = "2023-2024"
$ws.Range("E58").Value = 45192.54166666666
$ws.Range("F58").Value = "Horsens"
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = "Vendsyssel"
$ws.Range("I58").Value = 1
$ws.Range("J58").Value = 2.09
$ws.Range("K58").Value = "20/09/2023 00:13"
$ws.Range("L58").Value = 2.04
$ws.Range("M58").Value = "23/09/2023 12:31"
$ws.Range("N58").Value = 3.53
$ws.Range("O58").Value = "20/09/2023 00:13"
$ws.Range("P58").Value = 3.64
$ws.Range("Q58").Value = "23/09/2023 12:31"
$ws.Range("R58").Value = 3.21
$ws.Range("S58").Value = "20/09/2023 00:13"
$ws.Range("T58").Value = 3.55
$ws.Range("U58").Value = "23/09/2023 12:31"
$ws.Range("V58").Value = "https://www.betexplorer.com/football/denmark/1st-division/horsens-vendsyssel-ff/CUdyqBK9/"

$excel.CutCopyMode = 0
"done"